$wb = $excel.ActiveWorkbook

# --- studies sheet: add "PMID" column (H1) ---
$studies = $wb.Worksheets.Item("studies")
$studies.Range("H1").Value = "PMID"
$studies.Activate()
$studies.Range("H2").Select()

# --- counts sheet: add "notes" column (F1) ---
$counts = $wb.Worksheets.Item("counts")
$counts.Range("F1").Value = "notes"

# Make "counts" the active sheet / selection, as shown in diff (tabSelected moved from studies to counts)
$counts.Activate()
$counts.Range("F2").Select()
